# Updated Lectures 7, 8 and 9
#
# Slide 4  (08-Methods -> "public void print(String[] items)"):
#   merge the "public void " run with the "print(String" run.
# Slide 6  (CamelCase rule example):
#   merge "Camel" + "Case" runs into a single lower-camel-case run "camelCase".
# Slide 21 (exercise about a method multiplying 3 integers):
#   re-consolidate several adjacent runs that share identical formatting
#   (the text itself is unchanged, only the run/XML structure is simplified).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 4: "Пример:" -> "public void print(String[] items)"
# ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$para4 = $shp4.TextFrame.TextRange.Paragraphs(4, 1)
$para4.Characters(1, 24).Text = "public void print(String"

# ---------------------------------------------------------------
# Slide 6: "Трябва да се прилага правилото CamelCase." -> "...camelCase."
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$para6 = $shp6.TextFrame.TextRange.Paragraphs(4, 1)
$para6.Characters(32, 9).Text = "camelCase"

# ---------------------------------------------------------------
# Slide 21: "Дефинирайте метод, който приема 3 параметъра ..."
# ---------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$shp21 = $s21.Shapes.Item(2)
$para21 = $shp21.TextFrame.TextRange.Paragraphs(2, 1)

$para21.Characters(12, 6).Text = " метод"
$para21.Characters(25, 82).Text = " приема 3 параметъра - цели числа и ги умножава едно с друго. Извикайте функцията "
$para21.Characters(115, 28).Text = "с различни произволни числа."
